$wb = $excel.ActiveWorkbook

# Add the new worksheet "sheet2" after the existing Sheet1
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "sheet2"

# Row 1 - headers
$ws2.Range("A1").Value = "lower case"
$ws2.Range("B1").Value = "upper case"
$ws2.Range("C1").Value = "status"

# Row 2
$ws2.Range("A2").Value = "text"

# Row 3
$ws2.Range("A3").Value = "abc"
$ws2.Range("B3").Value = "ABC"
$ws2.Range("C3").Value = "Success"

# Row 4
$ws2.Range("A4").Value = "xyz"
$ws2.Range("B4").Value = "XYZ"
$ws2.Range("C4").Value = "Success"

# Row 5
$ws2.Range("A5").Value = "fgh"
$ws2.Range("B5").Value = "FGH"
$ws2.Range("C5").Value = "Success"
